# Update "想去人数" (F) and, where applicable, "最低票价" (G) counters
# across all 4 worksheets per the latest scrape snapshot.
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 167
$ws.Cells.Item(4, 6).Value = 2183
$ws.Cells.Item(5, 6).Value = 4335
$ws.Cells.Item(6, 6).Value = 574
$ws.Cells.Item(7, 6).Value = 1066
$ws.Cells.Item(8, 6).Value = 1337
$ws.Cells.Item(9, 6).Value = 667
$ws.Cells.Item(10, 6).Value = 378
$ws.Cells.Item(12, 6).Value = 405
$ws.Cells.Item(13, 6).Value = 671816
$ws.Cells.Item(13, 7).Value = 128
$ws.Cells.Item(14, 6).Value = 1670
$ws.Cells.Item(15, 6).Value = 573
$ws.Cells.Item(16, 6).Value = 1489
$ws.Cells.Item(17, 6).Value = 675
$ws.Cells.Item(19, 6).Value = 1304
$ws.Cells.Item(20, 6).Value = 2299
$ws.Cells.Item(21, 6).Value = 1162
$ws.Cells.Item(22, 6).Value = 2724
$ws.Cells.Item(23, 6).Value = 1577
$ws.Cells.Item(24, 6).Value = 872
$ws.Cells.Item(25, 6).Value = 1571
$ws.Cells.Item(28, 6).Value = 812
$ws.Cells.Item(29, 6).Value = 1103
$ws.Cells.Item(31, 6).Value = 90
$ws.Cells.Item(32, 6).Value = 2048
$ws.Cells.Item(33, 6).Value = 589
$ws.Cells.Item(34, 6).Value = 1331
$ws.Cells.Item(35, 6).Value = 2974
$ws.Cells.Item(36, 6).Value = 11
$ws.Cells.Item(37, 6).Value = 1154
$ws.Cells.Item(38, 6).Value = 46
$ws.Cells.Item(40, 6).Value = 2625
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(43, 6).Value = 3174
$ws.Cells.Item(44, 6).Value = 1027
$ws.Cells.Item(48, 6).Value = 673
$ws.Cells.Item(49, 6).Value = 23
$ws.Cells.Item(50, 6).Value = 23
# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(9, 6).Value = 114
$ws.Cells.Item(10, 6).Value = 494
$ws.Cells.Item(11, 6).Value = 144927
$ws.Cells.Item(12, 6).Value = 144927
$ws.Cells.Item(15, 6).Value = 23
$ws.Cells.Item(21, 6).Value = 424
$ws.Cells.Item(22, 6).Value = 170
$ws.Cells.Item(25, 6).Value = 91
$ws.Cells.Item(26, 6).Value = 625
$ws.Cells.Item(27, 6).Value = 91
$ws.Cells.Item(31, 6).Value = 371
$ws.Cells.Item(34, 6).Value = 71
$ws.Cells.Item(35, 6).Value = 71
$ws.Cells.Item(38, 6).Value = 225
$ws.Cells.Item(42, 6).Value = 89
# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 254
$ws.Cells.Item(7, 6).Value = 835
$ws.Cells.Item(8, 6).Value = 1219
$ws.Cells.Item(9, 6).Value = 645
$ws.Cells.Item(10, 6).Value = 1618
$ws.Cells.Item(11, 6).Value = 137
$ws.Cells.Item(12, 6).Value = 2018
# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 835
$ws.Cells.Item(3, 6).Value = 645
$ws.Cells.Item(5, 6).Value = 1618
$ws.Cells.Item(6, 6).Value = 167
$ws.Cells.Item(7, 6).Value = 2183
$ws.Cells.Item(8, 6).Value = 137
$ws.Cells.Item(9, 6).Value = 2018
$ws.Cells.Item(10, 6).Value = 4335
$ws.Cells.Item(11, 6).Value = 574
$ws.Cells.Item(12, 6).Value = 1337
$ws.Cells.Item(13, 6).Value = 667
$ws.Cells.Item(14, 6).Value = 378
$ws.Cells.Item(15, 6).Value = 405
$ws.Cells.Item(16, 6).Value = 671828
$ws.Cells.Item(16, 7).Value = 128
$ws.Cells.Item(17, 6).Value = 114
$ws.Cells.Item(18, 6).Value = 494
$ws.Cells.Item(19, 6).Value = 1670
$ws.Cells.Item(20, 6).Value = 144927
$ws.Cells.Item(21, 6).Value = 1489
$ws.Cells.Item(22, 6).Value = 675
$ws.Cells.Item(24, 6).Value = 1304
$ws.Cells.Item(25, 6).Value = 2299
$ws.Cells.Item(26, 6).Value = 1162
$ws.Cells.Item(27, 6).Value = 2724
$ws.Cells.Item(28, 6).Value = 1577
$ws.Cells.Item(29, 6).Value = 872
$ws.Cells.Item(30, 6).Value = 23
$ws.Cells.Item(31, 6).Value = 1571
$ws.Cells.Item(33, 6).Value = 170
$ws.Cells.Item(35, 6).Value = 812
$ws.Cells.Item(36, 6).Value = 1103
$ws.Cells.Item(37, 6).Value = 90
$ws.Cells.Item(38, 6).Value = 2048
$ws.Cells.Item(39, 6).Value = 1331
$ws.Cells.Item(40, 6).Value = 2974
$ws.Cells.Item(41, 6).Value = 11
$ws.Cells.Item(42, 6).Value = 1154
$ws.Cells.Item(43, 6).Value = 371
$ws.Cells.Item(45, 6).Value = 71
$ws.Cells.Item(46, 6).Value = 2625
$ws.Cells.Item(47, 6).Value = 216
$ws.Cells.Item(49, 6).Value = 3174
$ws.Cells.Item(50, 6).Value = 225
$ws.Cells.Item(51, 6).Value = 1027
$ws.Cells.Item(53, 6).Value = 673
$ws.Cells.Item(54, 6).Value = 23
